$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same bold/centered/bordered style
# used by the rest of row 1, instead of minting a brand-new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns for every data row (2-62): Wins=80, Losses=82, Ties=0
for ($row = 2; $row -le 62; $row++) {
    $ws.Cells.Item($row, 30).Value = 80
    $ws.Cells.Item($row, 31).Value = 82
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "done"
